$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.992.59'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.858.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.25'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3846'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08263'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -8.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.55'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.203'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.55'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.855.17'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.233'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.61'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06653'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.007'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.007.60'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.04'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.239'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.069.69'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.514'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.74'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.48'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.93'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1059'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.031'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.925'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.592'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.362'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02409'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06493'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2175'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6612'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.197'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.009'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.218'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.15'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6168'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.10'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.279'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.650'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.008'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.207'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '119.68'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.68'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.17%  '
